# Insert a new weekly price record at row 151 of the "Haba" price sheet.
# This shifts the existing rows 151-249 down to 152-250 (dimension becomes
# A1:R250), and the new row 151 is populated with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151, pushing rows 151..249 down to 152..250
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new data record
$ws.Range("A151").Value = 9
$ws.Range("B151").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C151").Value = "Metropolitana"
$ws.Range("D151").Value = 44777
$ws.Range("D151").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E151").Value = 13
$ws.Range("F151").Value = 100112026
$ws.Range("G151").Value = "Haba"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 52
$ws.Range("K151").Value = 16000
$ws.Range("L151").Value = 18000
$ws.Range("M151").Value = 17000
$ws.Range("N151").Value = "$/saco 25 kilos"
$ws.Range("O151").Value = "Región de Coquimbo"
$ws.Range("P151").Value = 680
$ws.Range("Q151").Value = 25
$ws.Range("R151").Value = "Hortaliza"
